$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Column widths (stored OOXML width = ColumnWidth + 5/6, so subtract 5/6 here
# to land on the exact integer widths used in the target file)
# ---------------------------------------------------------------------------
$offset = 5/6
$targetWidths = @(8,18,11,16,14,14,14,178,18,17,11,26,26,27,48,48,31,27,15,33,31,25,41,28)
for ($i = 0; $i -lt $targetWidths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $targetWidths[$i] - $offset
}

# ---------------------------------------------------------------------------
# Header row (row 1): center horizontally & vertically, wrap text.
# Build the final format once on a scratch cell (seeded from the existing
# header formatting so font/border carry over) and paste it onto the header
# range, so every header cell converges on a single shared style record.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats

$tmpl = $ws.Range("Z1")
$tmpl.WrapText = $true
$tmpl.HorizontalAlignment = -4108
$tmpl.VerticalAlignment = -4108

$tmpl.Copy()
$ws.Range("A1:X1").PasteSpecial(-4122)

# X1 additionally gets a yellow highlight fill on top of the header format
$tmpl.Interior.Color = 65535
$tmpl.Copy()
$ws.Range("X1").PasteSpecial(-4122)

$ws.Range("Z1").Clear()

# X1: updated label text
$ws.Range("X1").Value = "Status as of July 11, 2025"

# X2: yellow highlight fill (data cell under the updated header)
$ws.Range("X2").Interior.Color = 65535

# ---------------------------------------------------------------------------
# Freeze the header row (pane split after row 1) and keep A1 selected
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
